# Applies the "path_to_pbf" column addition + explanation text updates
# described in the commit "updated input excel and explanation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "input" sheet: insert a new column F ("path_to_pbf"), shifting the
#    existing shp_input_data / shp_unique_ID (and everything after them)
#    one column to the right.
# ---------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("input")
$wsInput.Columns("F").Insert()
$wsInput.Range("F1").Value = "path_to_pbf"

# Selection / view bookkeeping for the input sheet (no longer the active tab)
$wsInput.Range("D2").Select()

# ---------------------------------------------------------------------
# 2) "explanation" sheet: insert a new column F ("path_to_pbf") with its
#    explanation text, and rewrite the neighbouring explanation cells.
# ---------------------------------------------------------------------
$wsExpl = $wb.Worksheets.Item("explanation")
$wsExpl.Columns("F").Insert()
$wsExpl.Range("F1").Value = "path_to_pbf"
$wsExpl.Range("F2").Value = "When choosing network_source based on OSM dump. Provide name of *.pbf dump. No extension needed."

# E2 previously described the OSM option in general; now it specifically
# documents the "OSM online" option.
$wsExpl.Range("E2").Value = "When choosing network_source based on OSM online. Provide name of shapefile with region for OSM input <rel path is in config file>. No extension needed."

# G2 (formerly F2, the shapefile note with a red "consider removing" run)
# is replaced by plain, non-colored text.
$wsExpl.Range("G2").Value = "When choosing for network based on shapefile indicate shapefile for analysis.  Provide name of shapefile. No extension needed."

# ---------------------------------------------------------------------
# 3) Sheet view / active-tab bookkeeping: "explanation" becomes the
#    selected/active sheet instead of "input".
# ---------------------------------------------------------------------
$wsExpl.Range("G2").Select()
$wsExpl.Activate()

$wb.Save()
